$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Label" column (H) -----------------------------------
# Header cell, styled like the other header cells (bold + border).
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Label"

# Label values for every data row (Control rows = 0, MDD rows = 1)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1

# --- Refitted prediction/error values (re-run of the fit) -------------
$ws.Range("D3").Value = 0.411398788722149
$ws.Range("E3").Value = 0.411398788722149

$ws.Range("D4").Value = 0.4778773773823237
$ws.Range("E4").Value = 0.4778773773823237

$ws.Range("D6").Value = 0.4517313234496295
$ws.Range("E6").Value = 0.4517313234496295

$ws.Range("D7").Value = 0.4096532045775652
$ws.Range("E7").Value = 0.5903467954224348

$ws.Range("D11").Value = 0.4083546412490939
$ws.Range("E11").Value = 0.5916453587509061
$ws.Range("F11").Value = 0.6266096830368042
